$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rewrite -------------------------------------------------
# The sheet used to hold username/password/Tenant/Domain/POS (5 cols x 2 rows).
# It now holds just username/password across 3 rows (header + 2 credential rows).
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "lmenon"
$ws.Range("B2").Value = "Lavanya01"
$ws.Range("A3").Value = "Adas"
$ws.Range("B3").Value = "Welc0me03"

# Columns C:E (Tenant/Domain/POS) are no longer used - remove them entirely
# so the sheet dimension / spans shrink back down to A:B.
$ws.Range("C:E").Delete()

# Column B widened (no longer auto "best fit") to fit the new values.
$ws.Columns.Item(2).ColumnWidth = 16.71

# Bring the book window back to its normal working size.
$wb.Windows.Item(1).Width = 13335
$wb.Windows.Item(1).Height = 5175
